$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$data = @(
    @(45736.187303240738, 10, 6, 265, 453, 429, 478, 3437, 478, 2026, 208, 418, 30, 3885, 5002),
    @(45737.18372685185,  10, 6, 279, 453, 429, 478, 3437, 478, 2026, 208, 418, 30, 3935, 5042),
    @(45737.189606481479, 10, 6, 279, 453, 429, 478, 3437, 478, 2026, 208, 418, 30, 3935, 5042),
    @(45738.185208333336, 10, 6, 280, 454, 430, 478, 3437, 478, 2026, 208, 418, 30, 3941, 5064)
)

$startRow = 29
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
